# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E) and Pulse_Width (G)
# columns on each of the Step3_DataPts_* sheets to reflect the new
# zero_before_threshold behavior.

$wb = $excel.ActiveWorkbook

# Values to write for each sheet: row => @{ C = ...; E = ...; G = ... }
$updates = @{
    "Step3_DataPts_0.5" = @{
        2 = @{ C = 89; E = 0.001508937016090578; G = 40 }
        3 = @{ C = 88; E = 0.002599028523474619; G = 32 }
        4 = @{ C = 88; E = 0.001975553960472451; G = 34 }
        5 = @{ C = 88; E = 0.002925223205181372; G = 39 }
        6 = @{ C = 47; E = 0.0626444454130064;   G = 51 }
    }
    "Step3_DataPts_0.7" = @{
        2 = @{ C = 89; E = 0.001508937016090578; G = 58 }
        3 = @{ C = 88; E = 0.002599028523474619; G = 58 }
        4 = @{ C = 88; E = 0.001975553960472451; G = 57 }
        5 = @{ C = 88; E = 0.002925223205181372; G = 57 }
        6 = @{ C = 47; E = 0.0626444454130064;   G = 64 }
    }
    "Step3_DataPts_0.8" = @{
        2 = @{ C = 89; E = 0.001508937016090578; G = 65 }
        3 = @{ C = 88; E = 0.002599028523474619; G = 66 }
        4 = @{ C = 88; E = 0.001975553960472451; G = 65 }
        5 = @{ C = 88; E = 0.002925223205181372; G = 65 }
        6 = @{ C = 47; E = 0.0626444454130064;   G = 88 }
    }
    "Step3_DataPts_0.9" = @{
        2 = @{ C = 89; E = 0.001508937016090578; G = 77 }
        3 = @{ C = 88; E = 0.002599028523474619; G = 78 }
        4 = @{ C = 88; E = 0.001975553960472451; G = 76 }
        5 = @{ C = 88; E = 0.002925223205181372; G = 76 }
        6 = @{ C = 47; E = 0.0626444454130064;   G = 130 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $vals = $rows[$rowNum]
        $ws.Range("C$rowNum").Value = $vals.C
        $ws.Range("E$rowNum").Value = $vals.E
        $ws.Range("G$rowNum").Value = $vals.G
    }
}
